$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")
$ws.Range("C16").Value = 6
